# Generate Report for Handoff
# Updates the status of the 5166b202-... file from "Handed back: in sync with en-US"
# to "Ready for handoff" across the Overview, zh-cn and de-de sheets, and records the
# new "Latest Handoff Datetime" for that file in each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-03 12:42:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-03 12:42:30"
